$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F205").Value = 4935.333333333333
$ws.Range("G205").Value = 13655.33333333333
$ws.Range("H205").Value = 327589.6666666667
$ws.Range("F206").Value = 5007
$ws.Range("G206").Value = 13653.88888888889
$ws.Range("H206").Value = 324519.3333333333
$ws.Range("F207").Value = 4992.888888888889
$ws.Range("G207").Value = 13681.96296296296
$ws.Range("H207").Value = 323699.1111111111
$ws.Range("F208").Value = 4978.407407407407
$ws.Range("G208").Value = 13663.72839506173
$ws.Range("H208").Value = 325269.3703703704
$ws.Range("F209").Value = 4992.765432098765
$ws.Range("G209").Value = 13666.52674897119
$ws.Range("H209").Value = 324495.9382716049
$ws.Range("F210").Value = 4988.020576131687
$ws.Range("G210").Value = 13670.73936899863
$ws.Range("H210").Value = 324488.1399176955
$ws.Range("F211").Value = 4986.39780521262
$ws.Range("G211").Value = 13666.99817101052
$ws.Range("H211").Value = 324751.1495198903
$ws.Range("F212").Value = 4989.061271147691
$ws.Range("G212").Value = 13668.08809632678
$ws.Range("H212").Value = 324578.4092363969
$ws.Range("F213").Value = 4987.826550830666
$ws.Range("G213").Value = 13668.60854544531
$ws.Range("H213").Value = 324605.8995579942
$ws.Range("F214").Value = 4987.761875730325
$ws.Range("G214").Value = 13667.89827092754
$ws.Range("H214").Value = 324645.1527714271
$ws.Range("F215").Value = 4988.216565902894
$ws.Range("G215").Value = 13668.19830423321
$ws.Range("H215").Value = 324609.8205219394
$ws.Range("F216").Value = 4987.934997487962
$ws.Range("G216").Value = 13668.23504020202
$ws.Range("H216").Value = 324620.2909504536
$ws.Range("F217").Value = 4987.971146373727
$ws.Range("G217").Value = 13668.11053845425
$ws.Range("H217").Value = 324625.0880812734
$ws.Range("F218").Value = 4988.040903254861
$ws.Range("G218").Value = 13668.18129429649
$ws.Range("H218").Value = 324618.3998512221
$ws.Range("F219").Value = 4987.98234903885
$ws.Range("G219").Value = 13668.17562431759
$ws.Range("H219").Value = 324621.2596276497
$ws.Range("F220").Value = 4987.998132889145
$ws.Range("G220").Value = 13668.15581902278
$ws.Range("H220").Value = 324621.5825200484
$ws.Range("F221").Value = 4988.007128394285
$ws.Range("G221").Value = 13668.17091254562
$ws.Range("H221").Value = 324620.41399964
$ws.Range("F222").Value = 4987.995870107427
$ws.Range("G222").Value = 13668.16745196199
$ws.Range("H222").Value = 324621.085382446
$ws.Range("F223").Value = 4988.000377130285
$ws.Range("G223").Value = 13668.16472784346
$ws.Range("H223").Value = 324621.0273007115
$ws.Range("F224").Value = 4988.001125210666
$ws.Range("G224").Value = 13668.16769745036
$ws.Range("H224").Value = 324620.8422275992
$ws.Range("F225").Value = 4987.999124149459
$ws.Range("G225").Value = 13668.16662575194
$ws.Range("H225").Value = 324620.9849702522
$ws.Range("F226").Value = 4988.000208830137
$ws.Range("G226").Value = 13668.16635034859
$ws.Range("H226").Value = 324620.951499521
$ws.Range("F227").Value = 4988.000152730087
$ws.Range("G227").Value = 13668.16689118363
$ws.Range("H227").Value = 324620.9262324575
$ws.Range("F228").Value = 4987.999828569895
$ws.Range("G228").Value = 13668.16662242805
$ws.Range("H228").Value = 324620.9542340769
$ws.Range("F229").Value = 4988.000063376707
$ws.Range("G229").Value = 13668.16662132009
$ws.Range("H229").Value = 324620.9439886851
$ws.Range("F230").Value = 4988.00001489223
$ws.Range("G230").Value = 13668.16671164392
$ws.Range("H230").Value = 324620.9414850731
$ws.Range("F231").Value = 4987.999968946278
$ws.Range("G231").Value = 13668.16665179735
$ws.Range("H231").Value = 324620.9465692784
$ws.Range("F232").Value = 4988.000015738405
$ws.Range("G232").Value = 13668.16666158712
$ws.Range("H232").Value = 324620.9440143455
$ws.Range("F233").Value = 4987.999999858971
$ws.Range("G233").Value = 13668.16667500947
$ws.Range("H233").Value = 324620.944022899
$ws.Range("F234").Value = 4987.999994847884
$ws.Range("G234").Value = 13668.16666279798
$ws.Range("H234").Value = 324620.944868841
$ws.Range("F235").Value = 4988.000003481753
$ws.Range("G235").Value = 13668.16666646486
$ws.Range("H235").Value = 324620.9443020285
$ws.Range("F236").Value = 4987.999999396203
$ws.Range("G236").Value = 13668.16666809077
$ws.Range("H236").Value = 324620.9443979228
$ws.Range("F237").Value = 4987.999999241947
$ws.Range("G237").Value = 13668.16666578453
$ws.Range("H237").Value = 324620.9445229308
$ws.Range("F238").Value = 4988.000000706635
$ws.Range("G238").Value = 13668.16666678005
$ws.Range("H238").Value = 324620.9444076274
$ws.Range("F239").Value = 4987.999999781595
$ws.Range("G239").Value = 13668.16666688512
$ws.Range("H239").Value = 324620.944442827
$ws.Range("F240").Value = 4987.999999910059
$ws.Range("G240").Value = 13668.16666648324
$ws.Range("H240").Value = 324620.9444577951
$ws.Range("F241").Value = 4988.000000132763
$ws.Range("G241").Value = 13668.16666671614
$ws.Range("H241").Value = 324620.9444360831
$ws.Range("F242").Value = 4987.999999941472
$ws.Range("G242").Value = 13668.16666669483
$ws.Range("H242").Value = 324620.9444455684
$ws.Range("F243").Value = 4987.999999994765
$ws.Range("G243").Value = 13668.1666666314
$ws.Range("H243").Value = 324620.9444464822
$ws.Range("F244").Value = 4988.000000023
$ws.Range("G244").Value = 13668.16666668079
$ws.Range("H244").Value = 324620.9444427113
$ws.Range("F245").Value = 4987.999999986412
$ws.Range("G245").Value = 13668.16666666901
$ws.Range("H245").Value = 324620.9444449206
$ws.Range("F246").Value = 4988.000000001392
$ws.Range("G246").Value = 13668.1666666604
$ws.Range("H246").Value = 324620.9444447047
$ws.Range("F247").Value = 4988.000000003602
$ws.Range("G247").Value = 13668.16666667006
$ws.Range("H247").Value = 324620.9444441122
$ws.Range("F248").Value = 4987.999999997135
$ws.Range("G248").Value = 13668.16666666649
$ws.Range("H248").Value = 324620.9444445792
$ws.Range("F249").Value = 4988.000000000709
$ws.Range("G249").Value = 13668.16666666565
$ws.Range("H249").Value = 324620.9444444653
